# STM32G031_SOIC8_Dev_Board BOM update
# "3V3 LDO replaced, reverse polarity bug fixed"
#
# Updates the BOM worksheet (columns A-G: Comment, Description, Designator,
# Value, Quantity, DesignItemId, Footprint) to reflect several component
# substitutions:
#   - CN_8PIN connector -> CN_10PIN
#   - USB connector part 920-E52A2021S10100/USB307530A -> USB_B_AE
#   - Schottky diode BAT54C_R1_00001/SOT95P240X110-3N -> RB715UMTL/DAN217UMTL
#   - 3.3V LDO (IC2) NCP115ASN330T2G -> MIC5366-3.3YC5-TR (now 150mA, "LDO" desc)
#   - 1.8V LDO (IC3) MIC5365-1.8YD5-TR -> AP2120N-1.8TRG1 (now "LDO" desc,
#     footprint SOT95P230X110-3N)
#   - LEDs now specify current: Yellow -> Yellow/20mA, Green -> Green/20mA,
#     and new part numbers LTST-C171KSKT / LTST-C171KGKT
#   - New tactile Switch (S1, PTS815_SJM_250_SMTR_LFS) replaces the old
#     "Tactile switch" (SW1, EVP-BT3G4A000) row; the 0-ohm jumper (SB1) row
#     moves down to take its place

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: every text cell in this BOM table was originally stored with Excel's
# "quote prefix" (leading apostrophe) text style, which is what keeps
# part-number-like strings (e.g. "920-E52A2021S10100") from being
# reinterpreted as numbers/formulas/dates. Assigning a leading "'" on each
# new value preserves that same cell style (s="2") instead of resetting it
# to the default style COM would otherwise apply.

# Row 7: CN_8PIN -> CN_10PIN (connector for CN4, CN5)
$ws.Range("A7").Value = "'CN_10PIN"
$ws.Range("F7").Value = "'CN_10PIN"
$ws.Range("G7").Value = "'CN_10PIN"

# Row 8: USB connector part change
$ws.Range("A8").Value = "'USB_B_AE"
$ws.Range("F8").Value = "'USB_B_AE"
$ws.Range("G8").Value = "'USB_B_AE"

# Row 9: Schottky diode part change
$ws.Range("A9").Value = "'RB715UMTL"
$ws.Range("F9").Value = "'RB715UMTL"
$ws.Range("G9").Value = "'DAN217UMTL"

# Row 11: IC2, 3.3V LDO replaced
$ws.Range("A11").Value = "'MIC5366-3.3YC5-TR"
$ws.Range("B11").Value = "'LDO"
$ws.Range("D11").Value = "'3.3V/150mA"
$ws.Range("F11").Value = "'MIC5366-3.3YC5-TR"
$ws.Range("G11").Value = "'SOT65P210X110-5N"

# Row 12: IC3, 1.8V LDO replaced
$ws.Range("A12").Value = "'AP2120N-1.8TRG1"
$ws.Range("B12").Value = "'LDO"
$ws.Range("F12").Value = "'AP2120N-1.8TRG1"
$ws.Range("G12").Value = "'SOT95P230X110-3N"

# Row 13: LD1 Yellow LED part change
$ws.Range("A13").Value = "'LTST-C171KSKT"
$ws.Range("D13").Value = "'Yellow/20mA"
$ws.Range("F13").Value = "'LTST-C171KSKT"

# Row 14: LD2 Green LED part change
$ws.Range("A14").Value = "'LTST-C171KGKT"
$ws.Range("D14").Value = "'Green/20mA"
$ws.Range("F14").Value = "'LTST-C171KGKT"

# Row 16: new tactile Switch (S1) replaces the old jumper row content
$ws.Range("A16").Value = "'PTS815_SJM_250_SMTR_LFS"
$ws.Range("B16").Value = "'Switch"
$ws.Range("C16").Value = "'S1"
$ws.Range("F16").Value = "'PTS815_SJM_250_SMTR_LFS"
$ws.Range("G16").Value = "'PTS815SJM250SMTRLFS"

# Row 17: Jumper (0 ohm), SB1 -- now here instead of the old tactile switch row
$ws.Range("A17").Value = "'RC0805FR-070RL"
$ws.Range("B17").Value = "'Jumper (0 ohm)"
$ws.Range("C17").Value = "'SB1"
$ws.Range("F17").Value = "'RC0805FR-070RL"
$ws.Range("G17").Value = "'RESC2012X60N"
